$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Julián " + "Peker" (two runs, with proofErr spell-check markers around
# "Peker") become a single run "Julián Peker".
$d.Content.Find.Execute("Julián Peker", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Julián Peker", 2) | Out-Null

# --- Change 2 -----------------------------------------------------------
# "Resultado Esperado" cell: "...1m. Largo..." becomes
# "...1m. Altura máxima 2m. Largo..."
$d.Content.Find.Execute("1m. Largo", $false, $false, $false, $false, $false,
                         $true, 1, $false, "1m. Altura máxima 2m. Largo", 2) | Out-Null

# Same cell: "Largo mínima 1,8 m" becomes
# "Largo mínimo 1,8 m.  Largo máximo 2,70m."
$d.Content.Find.Execute("Largo mínima 1,8 m", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Largo mínimo 1,8 m.  Largo máximo 2,70m.", 2) | Out-Null

# --- Change 3 -----------------------------------------------------------
# "Resultado Obtenido" cell: " mínimo: 1,8 m y además se obtuvo el largo
# máximo 2,90 m" becomes " mínimo: 1,3 m y además se obtuvo el largo
# máximo 2,90 m y altura máxima 2,13m"
$d.Content.Find.Execute("mínimo: 1,8 m y además se obtuvo el largo máximo 2,90 m",
                         $false, $false, $false, $false, $false,
                         $true, 1, $false,
                         "mínimo: 1,3 m y además se obtuvo el largo máximo 2,90 m y altura máxima 2,13m", 2) | Out-Null
